$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2018) updates
$ws.Range("C2").Value = 24.37
$ws.Range("D2").Value = 6.72
$ws.Range("E2").Value = 60.96
$ws.Range("F2").Value = 45.65
$ws.Range("H2").Value = 89.97
$ws.Range("I2").Value = 85.28
$ws.Range("J2").Value = 37.81
$ws.Range("M2").Value = 42.88
$ws.Range("O2").Value = 70.77

# Row 3 (2019) updates
$ws.Range("H3").Value = 83.25
$ws.Range("I3").Value = 62.46

# Row 5 (2021) updates
$ws.Range("B5").Value = 0.45
$ws.Range("D5").Value = 14.69
$ws.Range("E5").Value = 94.90000000000001
$ws.Range("F5").Value = 78.86
$ws.Range("G5").Value = 94.75
$ws.Range("H5").Value = 94.90000000000001
$ws.Range("I5").Value = 44.08
$ws.Range("J5").Value = 48.28
$ws.Range("K5").Value = 44.83
$ws.Range("L5").Value = 16.04
$ws.Range("M5").Value = 80.06
$ws.Range("N5").Value = 33.88
$ws.Range("O5").Value = 53.82
$ws.Range("P5").Value = 93.09999999999999
